# "Analyses des performances après 2ème correction"
# Update the three score-tracking sheets with the latest measurements
# and leave the selection where the user last clicked on each tab.

$wb = $excel.ActiveWorkbook

# --- LightHouse - Portable --------------------------------------------
$wsPortable = $wb.Worksheets.Item("LightHouse - Portable")
$wsPortable.Range("E4").Value = 83
$wsPortable.Range("E4").Select()

# --- LightHouse - Bureau ------------------------------------------------
$wsBureau = $wb.Worksheets.Item("LightHouse - Bureau")
$wsBureau.Range("E4").Value = 92
$wsBureau.Range("E11").Select()

# --- GTmetrix - Bureau ---------------------------------------------------
$wsGTmetrix = $wb.Worksheets.Item("GTmetrix - Bureau")
$wsGTmetrix.Range("E5").Value = 96
$wsGTmetrix.Range("E5").Select()
